$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the firstname column (A2:A4) with new values
$ws.Range("A2").Value = "ABC"
$ws.Range("A3").Value = "DEF"
$ws.Range("A4").Value = "GHI"

# Move the active cell selection to H10 (matches diff: activeCell="H10" sqref="H10")
$ws.Range("H10").Select()
